$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(64, 8).Value = 10000
$ws.Cells.Item(64, 10).Value = 10000
$ws.Cells.Item(64, 12).Value = 10000
$ws.Cells.Item(64, 14).Value = -10496
$ws.Cells.Item(67, 8).Value = 10000
$ws.Cells.Item(67, 10).Value = 10000
$ws.Cells.Item(67, 12).Value = 10000
$ws.Cells.Item(67, 14).Value = -11716
$ws.Cells.Item(70, 8).Value = 12178
$ws.Cells.Item(70, 9).Value = 13549.875
$ws.Cells.Item(70, 11).Value = 40649.625
$ws.Cells.Item(70, 13).Value = -40379.625
$ws.Cells.Item(73, 8).Value = 12178
$ws.Cells.Item(73, 9).Value = 13549.875
$ws.Cells.Item(73, 11).Value = 40649.625
$ws.Cells.Item(73, 13).Value = -39713.625
$ws.Cells.Item(98, 8).Value = 1573.9678
$ws.Cells.Item(98, 9).Value = 1653.1305
$ws.Cells.Item(98, 11).Value = 1653.1305
$ws.Cells.Item(98, 13).Value = -155.1305
$ws.Cells.Item(107, 8).Value = 84373.664
$ws.Cells.Item(107, 9).Value = 1130.5
$ws.Cells.Item(107, 10).Value = 250860
$ws.Cells.Item(107, 11).Value = 1130.5
$ws.Cells.Item(107, 12).Value = 250860
$ws.Cells.Item(107, 13).Value = 789.5
$ws.Cells.Item(107, 14).Value = -254700
$ws.Cells.Item(122, 8).Value = 1573.9678
$ws.Cells.Item(122, 9).Value = 1653.1305
$ws.Cells.Item(122, 11).Value = 4959.3915
$ws.Cells.Item(122, 13).Value = -2509.3915
$ws.Cells.Item(138, 8).Value = 3046.6445
$ws.Cells.Item(138, 9).Value = 2421.8333
$ws.Cells.Item(138, 10).Value = 3273.8484
$ws.Cells.Item(138, 11).Value = 7265.499899999999
$ws.Cells.Item(138, 12).Value = 9821.5452
$ws.Cells.Item(138, 13).Value = -2125.499899999999
$ws.Cells.Item(138, 14).Value = -20101.5452

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 6244414.5
$ws.Cells.Item(32, 9).Value = 878458.5
$ws.Cells.Item(32, 11).Value = 878458.5
$ws.Cells.Item(32, 13).Value = -878171.5
$ws.Cells.Item(61, 8).Value = 1162.5
$ws.Cells.Item(61, 9).Value = 1186.3636
$ws.Cells.Item(61, 11).Value = 1186.3636
$ws.Cells.Item(61, 13).Value = -974.3635999999999
$ws.Cells.Item(88, 8).Value = 35824.11
$ws.Cells.Item(88, 9).Value = 1299.6666
$ws.Cells.Item(88, 10).Value = 53086.332
$ws.Cells.Item(88, 11).Value = 1299.6666
$ws.Cells.Item(88, 12).Value = 53086.332
$ws.Cells.Item(88, 13).Value = -893.6666
$ws.Cells.Item(88, 14).Value = -53898.332
$ws.Cells.Item(91, 8).Value = 35824.11
$ws.Cells.Item(91, 9).Value = 1299.6666
$ws.Cells.Item(91, 10).Value = 53086.332
$ws.Cells.Item(91, 11).Value = 1299.6666
$ws.Cells.Item(91, 12).Value = 53086.332
$ws.Cells.Item(91, 13).Value = 104.3334
$ws.Cells.Item(91, 14).Value = -55894.332
$ws.Cells.Item(94, 8).Value = 30330
$ws.Cells.Item(94, 10).Value = 30330
$ws.Cells.Item(94, 12).Value = 30330
$ws.Cells.Item(94, 14).Value = -32132
$ws.Cells.Item(95, 8).Value = 27604.5
$ws.Cells.Item(95, 9).Value = 25001
$ws.Cells.Item(95, 10).Value = 30208
$ws.Cells.Item(95, 11).Value = 25001
$ws.Cells.Item(95, 12).Value = 30208
$ws.Cells.Item(95, 13).Value = -22255
$ws.Cells.Item(95, 14).Value = -35700
$ws.Cells.Item(110, 8).Value = 1092.8695
$ws.Cells.Item(110, 10).Value = 2725
$ws.Cells.Item(110, 12).Value = 2725
$ws.Cells.Item(110, 14).Value = -6815
$ws.Cells.Item(132, 8).Value = 1789.4531
$ws.Cells.Item(132, 9).Value = 1641.0172
$ws.Cells.Item(132, 11).Value = 4923.0516
$ws.Cells.Item(132, 13).Value = -2393.0516
$ws.Cells.Item(136, 8).Value = 1162.5
$ws.Cells.Item(136, 9).Value = 1186.3636
$ws.Cells.Item(136, 11).Value = 3559.0908
$ws.Cells.Item(136, 13).Value = -1009.0908

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(64, 8).Value = 58052.332
$ws.Cells.Item(64, 10).Value = 69631.8
$ws.Cells.Item(64, 12).Value = 69631.8
$ws.Cells.Item(64, 14).Value = -70081.8
$ws.Cells.Item(67, 8).Value = 58052.332
$ws.Cells.Item(67, 10).Value = 69631.8
$ws.Cells.Item(67, 12).Value = 69631.8
$ws.Cells.Item(67, 14).Value = -71191.8
$ws.Cells.Item(105, 8).Value = 49306.332
$ws.Cells.Item(105, 9).Value = 3493.3333
$ws.Cells.Item(105, 10).Value = 72212.836
$ws.Cells.Item(105, 11).Value = 3493.3333
$ws.Cells.Item(105, 12).Value = 72212.836
$ws.Cells.Item(105, 13).Value = -1746.3333
$ws.Cells.Item(105, 14).Value = -75706.836
$ws.Cells.Item(134, 8).Value = 3109.717
$ws.Cells.Item(134, 9).Value = 2846.7297
$ws.Cells.Item(134, 10).Value = 3717.875
$ws.Cells.Item(134, 11).Value = 8540.1891
$ws.Cells.Item(134, 12).Value = 11153.625
$ws.Cells.Item(134, 13).Value = -6005.1891
$ws.Cells.Item(134, 14).Value = -16223.625

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 3841.8708
$ws.Cells.Item(31, 9).Value = 2360.8235
$ws.Cells.Item(31, 10).Value = 5640.2856
$ws.Cells.Item(31, 11).Value = 2360.8235
$ws.Cells.Item(31, 12).Value = 5640.2856
$ws.Cells.Item(31, 13).Value = -2065.8235
$ws.Cells.Item(31, 14).Value = -6230.2856
$ws.Cells.Item(34, 8).Value = 3841.8708
$ws.Cells.Item(34, 9).Value = 2360.8235
$ws.Cells.Item(34, 10).Value = 5640.2856
$ws.Cells.Item(34, 11).Value = 2360.8235
$ws.Cells.Item(34, 12).Value = 5640.2856
$ws.Cells.Item(34, 13).Value = -2158.8235
$ws.Cells.Item(34, 14).Value = -6044.2856
$ws.Cells.Item(62, 8).Value = 5814.5713
$ws.Cells.Item(62, 10).Value = 8000
$ws.Cells.Item(62, 12).Value = 8000
$ws.Cells.Item(62, 14).Value = -9248
$ws.Cells.Item(65, 8).Value = 5814.5713
$ws.Cells.Item(65, 10).Value = 8000
$ws.Cells.Item(65, 12).Value = 40000
$ws.Cells.Item(65, 14).Value = -46240
$ws.Cells.Item(112, 8).Value = 30000
$ws.Cells.Item(112, 10).Value = 30000
$ws.Cells.Item(112, 12).Value = 30000
$ws.Cells.Item(112, 14).Value = -32954

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 8).Value = 81.954544
$ws.Cells.Item(2, 9).Value = 47.785713
$ws.Cells.Item(2, 10).Value = 141.75
$ws.Cells.Item(2, 11).Value = 286.714278
$ws.Cells.Item(2, 12).Value = 850.5
$ws.Cells.Item(2, 13).Value = -173.714278
$ws.Cells.Item(2, 14).Value = -1076.5
$ws.Cells.Item(111, 8).Value = 5135.2
$ws.Cells.Item(111, 9).Value = 4419
$ws.Cells.Item(111, 10).Value = 8000
$ws.Cells.Item(111, 11).Value = 13257
$ws.Cells.Item(111, 12).Value = 24000
$ws.Cells.Item(111, 13).Value = -10190
$ws.Cells.Item(111, 14).Value = -30134
$ws.Cells.Item(122, 8).Value = 239.5
$ws.Cells.Item(122, 9).Value = 71.916664
$ws.Cells.Item(122, 10).Value = 490.875
$ws.Cells.Item(122, 11).Value = 647.2499759999999
$ws.Cells.Item(122, 12).Value = 4417.875
$ws.Cells.Item(122, 13).Value = 1802.750024
$ws.Cells.Item(122, 14).Value = -9317.875
$ws.Cells.Item(132, 8).Value = 1662.5
$ws.Cells.Item(132, 10).Value = 1702.1666
$ws.Cells.Item(132, 12).Value = 15319.4994
$ws.Cells.Item(132, 14).Value = -20379.4994

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 14).Value = 0
$ws.Cells.Item(43, 12).ClearContents()
$ws.Cells.Item(52, 8).Value = 36495
$ws.Cells.Item(52, 10).Value = 36495
$ws.Cells.Item(52, 12).Value = 36495
$ws.Cells.Item(52, 14).Value = -37013
$ws.Cells.Item(53, 8).Value = 62000
$ws.Cells.Item(53, 9).Value = 150000
$ws.Cells.Item(53, 10).Value = 18000
$ws.Cells.Item(53, 11).Value = 150000
$ws.Cells.Item(53, 12).Value = 18000
$ws.Cells.Item(53, 13).Value = -149369
$ws.Cells.Item(53, 14).Value = -19262
$ws.Cells.Item(80, 8).Value = 9431.799999999999
$ws.Cells.Item(80, 9).Value = 13288.111
$ws.Cells.Item(80, 11).Value = 13288.111
$ws.Cells.Item(80, 13).Value = -12290.111
$ws.Cells.Item(83, 8).Value = 9431.799999999999
$ws.Cells.Item(83, 9).Value = 13288.111
$ws.Cells.Item(83, 11).Value = 66440.55500000001
$ws.Cells.Item(83, 13).Value = -61448.55500000001

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 8).Value = 758332.5
$ws.Cells.Item(2, 10).Value = 758332.5
$ws.Cells.Item(2, 12).Value = 758332.5
$ws.Cells.Item(2, 14).Value = -758556.5
$ws.Cells.Item(68, 8).Value = 4534.4165
$ws.Cells.Item(68, 9).Value = 2426.625
$ws.Cells.Item(68, 11).Value = 2426.625
$ws.Cells.Item(68, 13).Value = -1677.625
$ws.Cells.Item(71, 8).Value = 4534.4165
$ws.Cells.Item(71, 9).Value = 2426.625
$ws.Cells.Item(71, 11).Value = 12133.125
$ws.Cells.Item(71, 13).Value = -8389.125
$ws.Cells.Item(82, 8).Value = 996.25
$ws.Cells.Item(82, 9).Value = 996.25
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 996.25
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).Value = -635.25
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(85, 8).Value = 996.25
$ws.Cells.Item(85, 9).Value = 996.25
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 996.25
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).Value = 251.75
$ws.Cells.Item(85, 13).ClearContents()

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(54, 8).Value = 30856.857
$ws.Cells.Item(54, 9).Value = 21500
$ws.Cells.Item(54, 10).Value = 43332.668
$ws.Cells.Item(54, 11).Value = 21500
$ws.Cells.Item(54, 12).Value = 43332.668
$ws.Cells.Item(54, 13).Value = -20980
$ws.Cells.Item(54, 14).Value = -44372.668
$ws.Cells.Item(74, 8).Value = 10515.667
$ws.Cells.Item(74, 10).Value = 10515.667
$ws.Cells.Item(74, 12).Value = 10515.667
$ws.Cells.Item(74, 14).Value = -12387.667
$ws.Cells.Item(77, 8).Value = 10515.667
$ws.Cells.Item(77, 10).Value = 10515.667
$ws.Cells.Item(77, 12).Value = 31547.001
$ws.Cells.Item(77, 14).Value = -40907.001
$ws.Cells.Item(81, 8).Value = 8908.666999999999
$ws.Cells.Item(81, 10).Value = 16666.334
$ws.Cells.Item(81, 12).Value = 33332.668
$ws.Cells.Item(81, 14).Value = -35454.668
$ws.Cells.Item(84, 8).Value = 8908.666999999999
$ws.Cells.Item(84, 10).Value = 16666.334
$ws.Cells.Item(84, 12).Value = 166663.34
$ws.Cells.Item(84, 14).Value = -177271.34
$ws.Cells.Item(107, 8).Value = 27810230
$ws.Cells.Item(107, 9).Value = 1209.375
$ws.Cells.Item(107, 10).Value = 50057450
$ws.Cells.Item(107, 11).Value = 3628.125
$ws.Cells.Item(107, 12).Value = 150172350
$ws.Cells.Item(107, 13).Value = -1708.125
$ws.Cells.Item(107, 14).Value = -150176190
$ws.Cells.Item(113, 8).Value = 2453.3333
$ws.Cells.Item(113, 9).Value = 1907
$ws.Cells.Item(113, 11).Value = 5721
$ws.Cells.Item(113, 13).Value = -3551
